$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the four "Jane McGonigal" rows (rows 11-14): blissful productivity,
# epic meaning, urgent optimism, social fabric. Everything below shifts up.
$ws.Range("A11:C14").EntireRow.Delete()

# Reset the selection to match the post-edit workbook state.
$ws.Range("B25").Select()
